$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the title cell (A1) to the new sheet title ("v1.3.2").
$ws.Range("A1").Value = "Microstate List"

# 2. Remove the two obsolete microstates (resonance structures / duplicate
#    geometric isomers SM20_micro001 and SM20_micro002), which occupied
#    rows 3-4. Deleting these rows shifts SM20_micro003/SM20_micro004 up
#    from rows 5-6 into rows 3-4.
$ws.Rows("3:4").Delete()

# 3. Remove the 2D-depiction pictures that belonged to the two deleted
#    microstates (the last two pictures anchored at the bottom of the
#    drawing canvas: "Picture 3" / "Picture 4").
$ws.Shapes.Item("Picture 4").Delete()
$ws.Shapes.Item("Picture 3").Delete()
